$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B4" = 0.6276404405018563
    "C4" = 0.06121979750756961
    "D4" = 0.5882352941176472
    "E4" = 0.5458515283842794
    "F4" = 0.6698113207547169
    "G4" = 0.6756756756756757
    "H4" = 0.5208333333333334
    "I4" = 0.6902884516993693
    "J4" = 0.7255747126436782
    "K4" = 46
    "L4" = 12
    "M4" = 23
    "N4" = 25
    "O4" = 0.7054367178970353
    "P4" = 0.02994035119622645

    "B5" = 0.6125956191124732
    "C5" = 0.005716248713015561
    "D5" = 0.6040268456375839
    "E5" = 0.7679180887372016
    "F5" = 0.4433962264150944
    "G5" = 0.4455445544554456
    "H5" = 0.9375
    "I5" = 0.4915980003721195
    "J5" = 0.4982040229885057
    "K5" = 2
    "L5" = 56
    "M5" = 3
    "N5" = 45
    "O5" = 0.512185689090451
    "P5" = 0.06000323185312067

    "B6" = 0.6108729882595384
    "C6" = 0.07140715576988323
    "D6" = 0.6526315789473683
    "E6" = 0.6485355648535565
    "F6" = 0.6886792452830188
    "G6" = 0.6595744680851063
    "H6" = 0.6458333333333334
    "I6" = 0.6267814364858217
    "J6" = 0.6968390804597702
    "K6" = 42
    "L6" = 16
    "M6" = 17
    "N6" = 31
    "O6" = 0.6550296044343662
    "P6" = 0.05511033388910008

    "B7" = 0.5776297060061557
    "C7" = 0.04822484832496659
    "D7" = 0.5176470588235293
    "E7" = 0.4803493449781659
    "F7" = 0.6132075471698113
    "G7" = 0.5945945945945946
    "H7" = 0.4583333333333333
    "I7" = 0.5178055413904471
    "J7" = 0.5998563218390806
    "K7" = 43
    "L7" = 15
    "M7" = 26
    "N7" = 22
    "O7" = 0.6183333333333334
    "P7" = 0.04527108237573843

    "B8" = 0.5398595012912438
    "C8" = 0.08100600269915967
    "O8" = 0.6723356009070296
    "P8" = 0.04971244811930053
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
